# Update Execution Time values in the "Add Majors Test" report sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "11389 ms"
$ws.Range("E3").Value = "6848 ms"
$ws.Range("E4").Value = "5747 ms"
$ws.Range("E5").Value = "5740 ms"
$ws.Range("E6").Value = "8916 ms"
$ws.Range("E7").Value = "7849 ms"
